$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The Price column (D) holds scraped values that look numeric (e.g. "1.001",
# "28.09") but must be written back as plain text, exactly matching the
# original inline-string cells (no locale/number coercion, no trailing-zero
# loss). Force Text format before assigning the value, then restore the
# cell style to Normal so the cell keeps the same (default) style it had
# before - only the stored string changes, just like in the source diff.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = '30.325.80'
$ws.Range("D3").Value = '1.933.43'
$ws.Range("D4").Value = '1.001'
$ws.Range("D5").Value = '0.7458'
$ws.Range("D6").Value = '250.15'
$ws.Range("D7").Value = '1.000'
$ws.Range("D8").Value = '28.09'
$ws.Range("D9").Value = '0.3219'
$ws.Range("D10").Value = '0.07142'
$ws.Range("D11").Value = '0.7903'
$ws.Range("D12").Value = '0.08007'
$ws.Range("D13").Value = '1.937.87'
$ws.Range("D14").Value = '5.389'
$ws.Range("D15").Value = '94.57'
$ws.Range("D16").Value = '14.59'
$ws.Range("D17").Value = '30.329.27'
$ws.Range("D18").Value = '252.58'
$ws.Range("D19").Value = '0.000008048'
$ws.Range("D20").Value = '5.806'
$ws.Range("D21").Value = '2.192.69'
$ws.Range("D22").Value = '1.000'
$ws.Range("D23").Value = '1.000'
$ws.Range("D24").Value = '6.822'
$ws.Range("D25").Value = '9.588'
$ws.Range("D26").Value = '164.38'
$ws.Range("D27").Value = '19.12'
$ws.Range("D28").Value = '2.315'
$ws.Range("D29").Value = '0.1332'
$ws.Range("D30").Value = '1.355'
$ws.Range("D31").Value = '1.535'
$ws.Range("D33").Value = '4.153'
$ws.Range("D34").Value = '0.05120'
$ws.Range("D35").Value = '1.290'
$ws.Range("D36").Value = '0.7489'
$ws.Range("D37").Value = '2.768'
$ws.Range("D38").Value = '0.01969'
$ws.Range("D39").Value = '2.803'
$ws.Range("D40").Value = '78.12'
$ws.Range("D41").Value = '6.416'
$ws.Range("D42").Value = '0.4524'
$ws.Range("D43").Value = '2.002'
$ws.Range("D44").Value = '0.8416'
$ws.Range("D45").Value = '1.000'
$ws.Range("D46").Value = '102.23'
$ws.Range("D47").Value = '7.587'
$ws.Range("D48").Value = '9.829'
$ws.Range("D49").Value = '987.32'
$ws.Range("D50").Value = '37.52'
$ws.Range("D51").Value = '0.1201'

foreach ($cell in $priceCells) {
    $ws.Range($cell).Style = "Normal"
}

# Coin name / link (B, C) and volume (E) columns are ordinary text strings
# (URLs, names, "  +1.23%  " style percentages) and can be assigned directly.
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("E3").Value = '  -0.41%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("E5").Value = '  +2.34%  '
$ws.Range("E6").Value = '  -0.42%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("B8").Value = 'Solana'
$ws.Range("C8").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("E8").Value = '  -2.67%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("E9").Value = '  -3.64%  '
$ws.Range("E10").Value = '  -2.67%  '
$ws.Range("E11").Value = '  -3.11%  '
$ws.Range("E12").Value = '  -1.27%  '
$ws.Range("E13").Value = '  -0.18%  '
$ws.Range("E14").Value = '  -1.83%  '
$ws.Range("E15").Value = '  -0.50%  '
$ws.Range("E16").Value = '  -2.91%  '
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("E19").Value = '  -3.17%  '
$ws.Range("E20").Value = '  -1.82%  '
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("E24").Value = '  -2.14%  '
$ws.Range("E25").Value = '  -2.50%  '
$ws.Range("E26").Value = '  +0.47%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("E27").Value = '  -1.59%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("E28").Value = '  -3.72%  '
$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("E29").Value = '  +0.65%  '
$ws.Range("E30").Value = '  +0.44%  '
$ws.Range("E31").Value = '  -2.64%  '
$ws.Range("E32").Value = '  -0.55%  '
$ws.Range("E33").Value = '  -1.78%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("E34").Value = '  -2.02%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("E35").Value = '  -1.66%  '
$ws.Range("E36").Value = '  -0.74%  '
$ws.Range("E37").Value = '  +1.00%  '
$ws.Range("E38").Value = '  -0.52%  '
$ws.Range("E39").Value = '  -1.62%  '
$ws.Range("E40").Value = '  -4.44%  '
$ws.Range("E41").Value = '  -2.33%  '
$ws.Range("E42").Value = '  -0.71%  '
$ws.Range("E43").Value = '  -2.21%  '
$ws.Range("E44").Value = '  -0.74%  '
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("E46").Value = '  -0.29%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("E47").Value = '  +1.08%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("E48").Value = '  -0.32%  '
$ws.Range("E49").Value = '  +12.02%  '
$ws.Range("E50").Value = '  +1.44%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("E51").Value = '  +5.03%  '
